# Updated symbol list on Mon Jan 30 13:41:42 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for the
# crypto rows on the worksheet, matching the latest scrape.
#
# Values are assigned with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr cells) instead of
# auto-converting numeric-looking strings/percentages into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.74"
$ws.Range("E2").Value = "'-2.21%"

$ws.Range("D3").Value = "'37.65"
$ws.Range("E3").Value = "'-4.49%"

$ws.Range("D4").Value = "'5.102"
$ws.Range("E4").Value = "'-0.19%"

$ws.Range("D5").Value = "'0.07870"
$ws.Range("E5").Value = "'-4.02%"

$ws.Range("D6").Value = "'1.989"
$ws.Range("E6").Value = "'1.21%"

$ws.Range("D7").Value = "'4.345"
$ws.Range("E7").Value = "'1.79%"

$ws.Range("D8").Value = "'8.240"
$ws.Range("E8").Value = "'-0.03%"

$ws.Range("D9").Value = "'3.131"
$ws.Range("E9").Value = "'-5.74%"

$ws.Range("D10").Value = "'0.9250"
$ws.Range("E10").Value = "'-0.65%"

$ws.Range("D11").Value = "'0.1276"
$ws.Range("E11").Value = "'-9.47%"

$ws.Range("D12").Value = "'0.1874"
$ws.Range("E12").Value = "'-4.95%"

$ws.Range("D13").Value = "'0.08814"
$ws.Range("E13").Value = "'-3.53%"

$ws.Range("D14").Value = "'0.03429"
$ws.Range("E14").Value = "'-3.11%"

$ws.Range("D15").Value = "'0.09765"
$ws.Range("E15").Value = "'-0.57%"

$ws.Range("D16").Value = "'0.001392"
$ws.Range("E16").Value = "'-0.71%"

$ws.Range("D17").Value = "'0.006047"
$ws.Range("E17").Value = "'1.71%"

$ws.Range("E18").Value = "'1,779.11%"

$ws.Range("D19").Value = "'3.576"
$ws.Range("E19").Value = "'-2.33%"

$ws.Range("D20").Value = "'0.3435"
$ws.Range("E20").Value = "'-0.81%"

$ws.Range("D21").Value = "'0.1284"
$ws.Range("E21").Value = "'-0.72%"

$ws.Range("D22").Value = "'5.010"
$ws.Range("E22").Value = "'2.39%"

$ws.Range("D23").Value = "'0.2498"
$ws.Range("E23").Value = "'2.11%"

$ws.Range("D24").Value = "'0.04327"
$ws.Range("E24").Value = "'0.17%"

$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'0.08%"

$ws.Range("D26").Value = "'0.004599"
$ws.Range("E26").Value = "'-4.00%"

$ws.Range("E27").Value = "'176.70%"

$ws.Range("D39").Value = "'0.02303"
$ws.Range("E39").Value = "'2.79%"

$ws.Range("D40").Value = "'0.05029"
$ws.Range("E40").Value = "'-4.66%"

$ws.Range("D41").Value = "'0.007529"
$ws.Range("E41").Value = "'0.14%"

$ws.Range("D42").Value = "'0.009847"
$ws.Range("E42").Value = "'-0.14%"

$ws.Range("D43").Value = "'0.1355"
$ws.Range("E43").Value = "'-1.82%"

$ws.Range("D44").Value = "'0.002094"
$ws.Range("E44").Value = "'-2.05%"

$ws.Range("D45").Value = "'0.008040"
$ws.Range("E45").Value = "'-17.89%"

$ws.Range("D46").Value = "'0.00006536"
$ws.Range("E46").Value = "'2.73%"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.44%"

$ws.Range("D48").Value = "'0.003006"
$ws.Range("E48").Value = "'8.80%"

$ws.Range("D49").Value = "'0.001204"
$ws.Range("E49").Value = "'0.44%"

$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.44%"

$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.44%"
